# Automatische test-sync: 2025-08-03 15:10:50
# Appends Testmail #18 ("Bestel je 200 stuks M8-bouten RVS voor Van Dijk?")
# as row 26 on the "Logs" sheet, extends the conditional-formatting ranges
# that cover the log table so they keep including the new row, and bumps
# the "Inkoop / Bestellingen" tally on the "Dashboard" sheet from 3 to 4.

$wb = $excel.ActiveWorkbook
$wsLogs = $wb.Worksheets.Item("Logs")

# --- Append the new log row (row 26) ---------------------------------
$wsLogs.Range("A26").Value = "Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$wsLogs.Range("B26").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C26").Value = "Testmail #18: Bestel je 200 stuks M8-bouten RVS voor Van Dijk?"
$wsLogs.Range("D26").Value = "Inkoop / Bestellingen"
$wsLogs.Range("E26").Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$wsLogs.Range("F26").Value = "2025-08-03 15:10:49"
$wsLogs.Range("G26").Value = "Ja"
$wsLogs.Range("H26").Value = "Ja"
$wsLogs.Range("I26").Value = "Nee"
$wsLogs.Range("J26").Value = "Nee"

# --- Extend conditional formatting so row 26 is covered too ----------
# (each of these columns has its own conditionalFormatting block in the
# sheet that previously applied to row 2 through row 25)
$wsLogs.Range("D2:D25").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("D2:D26"))
$wsLogs.Range("G2:G25").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("G2:G26"))
$wsLogs.Range("H2:H25").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("H2:H26"))
$wsLogs.Range("I2:I25").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("I2:I26"))
$wsLogs.Range("J2:J25").FormatConditions.Item(1).ModifyAppliesToRange($wsLogs.Range("J2:J26"))

# --- Update the Dashboard summary count for "Inkoop / Bestellingen" --
$wsDash = $wb.Worksheets.Item("Dashboard")
$wsDash.Range("B5").Value = 4
